# Update the cryptocurrency price ("Price", column D) and 1-hour volume
# change percentage ("Volume(1h)", column E) figures for each coin row
# (rows 2-51) on the active sheet, reflecting the latest scrape refresh.
#
# D34 is forced back to literal text (apostrophe-prefixed, then restyled
# to "Normal") because "23.50" would otherwise be auto-coerced to the
# number 23.5 and display its trailing zero stripped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.284.37'
$ws.Range("E2").Value = '  +0.59%  '
$ws.Range("D3").Value = '3.483.36'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '593.82'
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '178.22'
$ws.Range("E6").Value = '  +3.65%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  +1.59%  '
$ws.Range("D9").Value = '3.486.31'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("E10").Value = '  +4.69%  '
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("D12").Value = '0.435'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '4.084.53'
$ws.Range("E13").Value = '  -0.35%  '
$ws.Range("D14").Value = '31.91'
$ws.Range("E14").Value = '  +9.62%  '
$ws.Range("D15").Value = '0.136'
$ws.Range("E15").Value = '  +1.04%  '
$ws.Range("D16").Value = '67.302.93'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("E17").Value = '  -0.98%  '
$ws.Range("D18").Value = '3.481.84'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '14.27'
$ws.Range("E20").Value = '  +1.61%  '
$ws.Range("D21").Value = '388.49'
$ws.Range("E21").Value = '  -1.58%  '
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("D23").Value = '73.85'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '0.537'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("E29").Value = '  -3.91%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").Value = '2.06'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = "'23.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("D35").Value = '7.36'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = '164.61'
$ws.Range("E38").Value = '  +1.07%  '
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").Value = '2.75'
$ws.Range("E41").Value = '  +7.99%  '
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("D43").Value = '4.65'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = '2.832.42'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '27.03'
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").Value = '26.13'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("E47").Value = '  -2.51%  '
$ws.Range("D48").Value = '41.61'
$ws.Range("E48").Value = '  -2.80%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").Value = '334.73'
$ws.Range("E50").Value = '  -0.35%  '
$ws.Range("E51").Value = '  -2.30%  '
